# "Added Old Slides + Updated Syllabus" - reorders/updates the session
# topics in column C (rows 8-18), adds a new highlighted "library" session,
# splits the "Quantitative Methoden" block into separate sessions, marks
# two new placeholder rows (27, 28) with matching formatting/conditional
# formatting, and leaves a new selection on C21:C24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign the session titles in column C (rows 8-18) -------------------
# Rows 4,5,6,7,19 keep their original text - only the ones below actually
# change content (some are simply moved to a new row, some are brand new).
$ws.Range("C8").Value  = "(Bibliothek + Asynchron) Wahlen und Wähler in Deutschland und der Welt "
$ws.Range("C9").Value  = "(Asynchron) Quantitative Methoden I - Einführung in R "
$ws.Range("C10").Value = "(Präsenz) Wer wählt populistische Parteien? + Wissenschaftliches Lesen I"
$ws.Range("C11").Value = "(Asynchron) Quantitative Methoden II – Daten reinigen und visualisieren "
$ws.Range("C12").Value = "(Präsenz) Wissenschaftliches Lesen II: Literaturreviews und Recherche "
$ws.Range("C13").Value = "(Online) Sprechstunde für Hausarbeiten 1"
$ws.Range("C14").Value = "(Präsenz) Quantitative Methoden III – Logik der quantitativen Methoden "
$ws.Range("C15").Value = "(Asynchron) Quantitative Methoden IV – Anwendung 1"
$ws.Range("C16").Value = "(Präsenz) Quantitative Methoden V – Anwendung 2"
$ws.Range("C17").Value = "(Präsenz) Qualitative Methoden "
$ws.Range("C18").Value = "(Online) Sprechstunde für Hausarbeiten 2"

# --- Highlight the new "library" session (C8) in green ---------------------
$c8 = $ws.Range("C8")
$c8.Font.Name = "Times New Roman"
$c8.Font.Size = 12
$c8.Interior.Color = 5287936
$c8.HorizontalAlignment = -4131
$c8.VerticalAlignment = -4108

# --- C14 moves from the plain style onto the standard session style --------
$c14 = $ws.Range("C14")
$c14.Font.Name = "Times New Roman"
$c14.Font.Size = 12
$c14.HorizontalAlignment = -4131
$c14.VerticalAlignment = -4108

# --- New placeholder rows --------------------------------------------------
$c27 = $ws.Range("C27")
$c27.Font.Name = "Calibri"
$c27.Font.Size = 11
$c27.Font.Bold = $true

$c28 = $ws.Range("C28")
$c28.Font.Name = "Times New Roman"
$c28.Font.Size = 12
$c28.HorizontalAlignment = -4131
$c28.VerticalAlignment = -4108
$ws.Rows.Item(28).RowHeight = 15.75

# --- Conditional formatting: bump priorities of the existing rules, then ---
# --- add two new "(Präsenz)" highlight rules matching the new rows --------
$r1 = $ws.Range("C4:C19")
$f1 = $r1.FormatConditions.Item(1)
$f1.Priority = 3

$r2 = $ws.Range("C30")
$f2 = $r2.FormatConditions.Item(1)
$f2.Priority = 5

$rngC28 = $ws.Range("C28")
$fcC28 = $rngC28.FormatConditions.Add(2, 0, 'ISNUMBER(SEARCH("(Präsenz)", C28:C43))')
$fcC28.Interior.Color = 5287936
$fcC28.Priority = 1

$rngC24 = $ws.Range("C24")
$fcC24 = $rngC24.FormatConditions.Add(2, 0, 'ISNUMBER(SEARCH("(Präsenz)", C15:C29))')
$fcC24.Interior.Color = 5287936
$fcC24.Priority = 7

# --- Final selection ---------------------------------------------------
$ws.Range("C21:C24").Select()
